$wb = $excel.ActiveWorkbook

# --- Update the conversion note on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = [char]0x2705 + " 1000 Bs = 6.98 = 27308.1 pesos"
$newLine1 = [char]0x2705 + " 1000 Bs = 7.03 = 27558.37 pesos"
$oldLine2 = [char]0x2705 + " 27308.1 pesos = 6.95 = 966.55 Bs"
$newLine2 = [char]0x2705 + " 27558.37 pesos = 7.0 = 973.48 Bs"

[string]$text = $ws1.Range("A1").Value2
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$ws1.Range("A1").Value = $text

# --- Update the rate cells on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value2 = 142.2
$ws2.Range("O10").Value2 = 3918.8
$ws2.Range("N12").Value2 = 3934.97
$ws2.Range("O12").Value2 = 139
